$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new values look like numbers (e.g. "1.005") need to be
# forced to Text format first, otherwise Excel auto-converts them to
# numeric values and strips meaningful trailing/leading zeros. ---
$textForceRefs = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$newValues = @{
    "D4" = "1.005"
    "D5" = "244.91"
    "D6" = "0.6313"
    "D7" = "1.006"
    "D8" = "0.07510"
    "D9" = "0.2937"
    "D10" = "23.16"
    "D11" = "0.07751"
    "D13" = "4.998"
    "D14" = "0.6705"
    "D15" = "83.14"
    "D16" = "0.000009356"
    "D17" = "6.055"
    "D19" = "12.64"
    "D20" = "223.98"
    "D21" = "1.007"
    "D22" = "7.153"
    "D23" = "1.007"
    "D24" = "160.93"
    "D25" = "0.1403"
    "D26" = "8.536"
    "D27" = "17.97"
    "D28" = "1.511"
    "D29" = "4.166"
    "D30" = "0.05572"
    "D31" = "4.085"
    "D32" = "1.209"
    "D33" = "0.7547"
    "D34" = "1.859"
    "D35" = "1.140"
    "D36" = "2.618"
    "D38" = "0.01792"
    "D39" = "2.757"
    "D40" = "6.580"
    "D41" = "0.8978"
    "D42" = "1.007"
    "D43" = "102.19"
    "D45" = "65.95"
    "D46" = "0.00000000125"
    "D47" = "0.07805"
    "D48" = "0.5107"
    "D49" = "0.4072"
    "D50" = "9.081"
    "D51" = "0.05818"
    "D2" = "29.192.70"
    "E2" = "  +0.24%  "
    "D3" = "1.839.33"
    "E3" = "  +0.05%  "
    "E4" = "  +0.27%  "
    "E5" = "  +0.92%  "
    "E6" = "  +0.93%  "
    "E7" = "  +0.28%  "
    "E8" = "  -1.55%  "
    "E9" = "  +0.46%  "
    "E10" = "  +1.80%  "
    "E11" = "  -0.10%  "
    "D12" = "1.837.65"
    "E12" = "  -0.13%  "
    "E13" = "  +0.69%  "
    "E14" = "  +0.47%  "
    "E15" = "  +0.20%  "
    "E16" = "  -7.67%  "
    "E17" = "  +0.51%  "
    "D18" = "29.206.89"
    "E18" = "  +0.20%  "
    "E19" = "  +2.15%  "
    "E20" = "  -0.88%  "
    "E21" = "  +0.42%  "
    "E22" = "  -0.96%  "
    "E23" = "  +0.27%  "
    "E24" = "  +1.38%  "
    "E25" = "  +2.61%  "
    "E26" = "  +0.91%  "
    "E27" = "  +0.19%  "
    "E28" = "  +1.39%  "
    "E29" = "  +2.13%  "
    "E30" = "  +6.95%  "
    "E31" = "  +1.11%  "
    "E32" = "  +0.41%  "
    "E33" = "  +2.01%  "
    "E34" = "  +0.38%  "
    "E35" = "  -0.59%  "
    "E36" = "  -3.30%  "
    "D37" = "1.237.73"
    "E37" = "  -2.13%  "
    "B38" = "VeChain"
    "C38" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "E38" = "  +0.41%  "
    "B39" = "MXToken"
    "C39" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "E39" = "  -0.46%  "
    "E40" = "  +3.93%  "
    "E41" = "  -0.12%  "
    "E42" = "  +0.28%  "
    "E43" = "  +0.44%  "
    "D44" = "1.985.53"
    "E44" = "  +0.36%  "
    "B45" = "Aave"
    "C45" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "E45" = "  +2.82%  "
    "B46" = "BabyDogeCoin"
    "C46" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "E46" = "  +0.22%  "
    "E47" = "  +14.27%  "
    "E48" = "  -0.32%  "
    "E49" = "  +1.11%  "
    "E50" = "  +2.34%  "
    "E51" = "  +1.04%  "
}
foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value = $newValues[$ref]
}